# Update latest output (run 193)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule" (1st sheet)
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item(1)

# Update existing row 2 (Cost / Unit Cost recalculated)
$schedule.Cells.Item(2, 5).Value = 505.13873475
$schedule.Cells.Item(2, 6).Value = 12.14859871933622

# Add the new schedule row 3, matching the date/time format used in row 2
$schedule.Range("A3:B3").NumberFormat = $schedule.Range("A2:B2").NumberFormat

$schedule.Cells.Item(3, 1).Value = 46046.3125
$schedule.Cells.Item(3, 2).Value = 46046.77083333334
$schedule.Cells.Item(3, 3).Value = 11
$schedule.Cells.Item(3, 4).Value = 41.58
$schedule.Cells.Item(3, 5).Value = -245.24632275
$schedule.Cells.Item(3, 6).Value = -5.898179960317462

# ---------------------------------------------------------------------------
# Sheet "Detailed" (2nd sheet)
# ---------------------------------------------------------------------------
$detailed = $wb.Worksheets.Item(2)

# Updated Price values for existing rows 27-47
$priceUpdates = @(
    @(27, 33.78973),
    @(28, 34.45564),
    @(29, 36.06),
    @(30, 0.01698),
    @(31, 11.19473),
    @(32, 0.02584),
    @(33, 0.51),
    @(34, 0.51),
    @(35, 9.20463),
    @(36, 22.07),
    @(37, 38.64756),
    @(38, 48.43349),
    @(39, 57.06),
    @(40, 63.33846),
    @(41, 66.43913999999999),
    @(42, 68.9002),
    @(43, 63.53656),
    @(44, 73.2),
    @(45, 65),
    @(46, 59.96931),
    @(47, 57.31)
)

foreach ($item in $priceUpdates) {
    $rowNum = $item[0]
    $price = $item[1]
    $detailed.Cells.Item($rowNum, 2).Value = $price
}

# Rows 29-32 flip from "forecast" to "historical"
foreach ($rowNum in 29..32) {
    $detailed.Cells.Item($rowNum, 3).Value = "historical"
}

# New rows 50-97 appended after the existing data (row 49 was the last one)
$newRows = @(
    @(50, 46046, 57.06, "forecast", 46046, "OFF"),
    @(51, 46046.02083333334, 59.66067, "forecast", 46046, "OFF"),
    @(52, 46046.04166666666, 57.06, "forecast", 46046, "OFF"),
    @(53, 46046.0625, 57.06, "forecast", 46046, "OFF"),
    @(54, 46046.08333333334, 57.06, "forecast", 46046, "OFF"),
    @(55, 46046.10416666666, 57.06, "forecast", 46046, "OFF"),
    @(56, 46046.125, 57.06, "forecast", 46046, "OFF"),
    @(57, 46046.14583333334, 57.06, "forecast", 46046, "OFF"),
    @(58, 46046.16666666666, 63.25033, "forecast", 46046, "OFF"),
    @(59, 46046.1875, 73.2, "forecast", 46046, "OFF"),
    @(60, 46046.20833333334, 73.2, "forecast", 46046, "OFF"),
    @(61, 46046.22916666666, 73.2, "forecast", 46046, "OFF"),
    @(62, 46046.25, 66.12326, "forecast", 46046, "OFF"),
    @(63, 46046.27083333334, 59.56115, "forecast", 46046, "OFF"),
    @(64, 46046.29166666666, 36.06, "forecast", 46046, "OFF"),
    @(65, 46046.3125, 0.7, "forecast", 46046, "ON"),
    @(66, 46046.33333333334, -5.78628, "forecast", 46046, "ON"),
    @(67, 46046.35416666666, -5.74313, "forecast", 46046, "ON"),
    @(68, 46046.375, -7.2053, "forecast", 46046, "ON"),
    @(69, 46046.39583333334, -12.01, "forecast", 46046, "ON"),
    @(70, 46046.41666666666, -14, "forecast", 46046, "ON"),
    @(71, 46046.4375, -14, "forecast", 46046, "ON"),
    @(72, 46046.45833333334, -20, "forecast", 46046, "ON"),
    @(73, 46046.47916666666, -23.69391, "forecast", 46046, "ON"),
    @(74, 46046.5, -27, "forecast", 46046, "ON"),
    @(75, 46046.52083333334, -27, "forecast", 46046, "ON"),
    @(76, 46046.54166666666, -30.05223, "forecast", 46046, "ON"),
    @(77, 46046.5625, -29.5, "forecast", 46046, "ON"),
    @(78, 46046.58333333334, -27, "forecast", 46046, "ON"),
    @(79, 46046.60416666666, -13.5, "forecast", 46046, "ON"),
    @(80, 46046.625, -7.20607, "forecast", 46046, "ON"),
    @(81, 46046.64583333334, -7.47463, "forecast", 46046, "ON"),
    @(82, 46046.66666666666, -6, "forecast", 46046, "ON"),
    @(83, 46046.6875, -6, "forecast", 46046, "ON"),
    @(84, 46046.70833333334, -5.50985, "forecast", 46046, "ON"),
    @(85, 46046.72916666666, 8.63959, "forecast", 46046, "ON"),
    @(86, 46046.75, 27.80712, "forecast", 46046, "ON"),
    @(87, 46046.77083333334, 57.68272, "forecast", 46046, "OFF"),
    @(88, 46046.79166666666, 74.98204, "forecast", 46046, "OFF"),
    @(89, 46046.8125, 91.00466, "forecast", 46046, "OFF"),
    @(90, 46046.83333333334, 90.84325, "forecast", 46046, "OFF"),
    @(91, 46046.85416666666, 105, "forecast", 46046, "OFF"),
    @(92, 46046.875, 95.04008, "forecast", 46046, "OFF"),
    @(93, 46046.89583333334, 93.89266000000001, "forecast", 46046, "OFF"),
    @(94, 46046.91666666666, 83.37803, "forecast", 46046, "OFF"),
    @(95, 46046.9375, 108.89, "forecast", 46046, "OFF"),
    @(96, 46046.95833333334, 105, "forecast", 46046, "OFF"),
    @(97, 46046.97916666666, 92.51152, "forecast", 46046, "OFF")
)

# Apply the same date/time formats used on the existing rows to the new ones
$detailed.Range("A50:A97").NumberFormat = $detailed.Range("A49").NumberFormat
$detailed.Range("D50:D97").NumberFormat = $detailed.Range("D49").NumberFormat

foreach ($row in $newRows) {
    $rowNum = $row[0]
    $detailed.Cells.Item($rowNum, 1).Value = $row[1]
    $detailed.Cells.Item($rowNum, 2).Value = $row[2]
    $detailed.Cells.Item($rowNum, 3).Value = $row[3]
    $detailed.Cells.Item($rowNum, 4).Value = $row[4]
    $detailed.Cells.Item($rowNum, 5).Value = $row[5]
}
